# Tasks.xlsx edit: add a "description" column, add 3 new task rows, assign
# "Pooriya" as the doer of the Header/Footer-related tasks, and register the
# LOCAL_MYSQL_DATE_FORMAT hidden defined name (MySQL for Excel add-in
# leftover).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Hidden workbook-level defined name (MySQL add-in leftover artifact).
# ---------------------------------------------------------------------------
$mysqlDateFormula = '=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&" "&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)'
$definedName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $mysqlDateFormula)
$definedName.Visible = $false

# ---------------------------------------------------------------------------
# 2. Insert three new rows right above the old last row ("تعریف تیم...")
#    so it moves from row 12 down to row 15.
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).RowHeight = 18.6
$ws.Rows.Item(13).RowHeight = 18.6
$ws.Rows.Item(14).RowHeight = 18.6

# ---------------------------------------------------------------------------
# 3. New column C ("توضیحات" / description), with a wider custom width.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "توضیحات"
$ws.Columns.Item(3).ColumnWidth = 80.1

# ---------------------------------------------------------------------------
# 4. Fill in the task text / doer / description for every data row.
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = "پوریا"

$ws.Range("A12").Value = "صفحه ی ثبت نام و ورود به سیستم"
$ws.Range("B12").Value = "پوریا"
$ws.Range("C12").Value = "در Header وارد شده است. هر سه نوع را پذیرش می کند. تنها کافی است فیلدهای ثبت نام بازبینی و دقیقتر شوند."

$ws.Range("A13").Value = "ساخت Footer"
$ws.Range("B13").Value = "پوریا"

$ws.Range("A14").Value = "پیدا کردن عکس مناسب با اندازه ی 1920*500 برای اسلایدر"

# ---------------------------------------------------------------------------
# 5. Formatting.
#    - Header row (A1:C1) all share the bold, centered header look.
#    - The "doer" column (B) is centered for every data row.
#    - The new description column (C) uses the same Persian font as the
#      rest of the table; row 12's description cell keeps the plain
#      (non-centered) look, like column A.
# ---------------------------------------------------------------------------
$ws.Range("C1").Font.Name = "B Nazanin"
$ws.Range("C1").Font.Size = 12
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108

$ws.Range("B2:B15").HorizontalAlignment = -4108

$ws.Range("C2:C11").Font.Name = "B Nazanin"
$ws.Range("C2:C11").Font.Size = 12
$ws.Range("C13:C15").Font.Name = "B Nazanin"
$ws.Range("C13:C15").Font.Size = 12

$ws.Range("C12").Font.Name = "B Nazanin"
$ws.Range("C12").Font.Size = 12
$ws.Range("C13").Font.Name = "B Nazanin"
$ws.Range("C13").Font.Size = 12
$ws.Range("C14").Font.Name = "B Nazanin"
$ws.Range("C14").Font.Size = 12

# ---------------------------------------------------------------------------
# 6. Resize the Table1 list object to the new A1:C15 range.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C15"))
$ws.Range("C1").Value = "توضیحات"

# ---------------------------------------------------------------------------
# 7. Selection cosmetics to match the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("B11").Select()
